$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: drop the explicit "no underline" direct formatting
#    (<w:u w:val="none"/>) from the run properties, keeping the Arial /
#    black / 28-half-point (14pt) formatting intact.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleText = $d.Range($titleRange.Start, $titleRange.End - 1).Text

$textOnly = $d.Range($titleRange.Start, $titleRange.End - 1)
$textOnly.Delete()

$refreshedPara = $d.Paragraphs.Item(1)
$refreshedPara.Range.InsertBefore($titleText)

$newTitlePara = $d.Paragraphs.Item(1)
$newTitleText = $d.Range($newTitlePara.Range.Start, $newTitlePara.Range.End - 1)
$newTitleText.Font.Name = "Arial"
$newTitleText.Font.Color = 0
$newTitleText.Font.Size = 14

# ---------------------------------------------------------------------------
# 2) Header row (row 1) of the table: insert a blank paragraph before the
#    existing content paragraph in each of the 4 cells.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
for ($i = 1; $i -le 4; $i++) {
    $cell = $t.Cell(1, $i)
    $cell.Range.InsertParagraphBefore()
    $blankPara = $cell.Range.Paragraphs.Item(1)
    $blankPara.Style = "Normal"
    $blankPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Append a brand-new, entirely empty row (4 cells, each just a blank
#    paragraph) at the bottom of the table.
# ---------------------------------------------------------------------------
$newRow = $t.Rows.Add()
$newRowIndex = $t.Rows.Count
for ($i = 1; $i -le 4; $i++) {
    $cell = $t.Cell($newRowIndex, $i)
    $blankPara = $cell.Range.Paragraphs.Item(1)
    $blankPara.Style = "Normal"
    $blankPara.Range.Delete()
}

Write-Host "Edit complete."
